$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the date serial value in A2:A9 from 45779 to 45780 (one day later)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = 45780
}
